$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FRED Graph")

# Revised GDP observations (2001-2019) for existing rows 12-30
$ws.Range("B12").Value = 234693.57
$ws.Range("B13").Value = 237735.492
$ws.Range("B14").Value = 245407.958
$ws.Range("B15").Value = 259116.344
$ws.Range("B16").Value = 272211.311
$ws.Range("B17").Value = 287146.989
$ws.Range("B18").Value = 303949.312
$ws.Range("B19").Value = 311634.603
$ws.Range("B20").Value = 312593.395
$ws.Range("B21").Value = 330424.307
$ws.Range("B22").Value = 341801.133
$ws.Range("B23").Value = 357088.355
$ws.Range("B24").Value = 366244.863
$ws.Range("B25").Value = 379924.795
$ws.Range("B26").Value = 403931.138
$ws.Range("B27").Value = 418289.875
$ws.Range("B28").Value = 434282.612
$ws.Range("B29").Value = 462229.528
$ws.Range("B30").Value = 488174.288

# New observation for 2020-01-01, matching the existing date/number formats
$ws.Range("A31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("A31").Value = 43831
$ws.Range("B31").NumberFormat = "0.000"
$ws.Range("B31").Value = 480307.073

# Restore the sheet's selection state: the two data columns (A:B) selected
$ws.Range("A1:B1048576").Select()
